$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = [double]"0.08850642943357422"

$ws.Range("B3").Value = [double]"0.0005939510681812405"
$ws.Range("C3").Value = [double]"0.0005464916938715664"
$ws.Range("D3").Value = [double]"2.545066805924304"
$ws.Range("E3").Value = [double]"0.2144713576259772"
$ws.Range("F3").Value = [double]"-0.0004771560516121309"
$ws.Range("G3").Value = [double]"0.001665058187974612"
$ws.Range("H3").Value = [double]"0.08910038050175545"

$ws.Range("B4").Value = [double]"0.009935773338771449"
$ws.Range("C4").Value = [double]"0.0007124288880028367"
$ws.Range("D4").Value = [double]"13.9138704988514"
$ws.Range("E4").Value = [double]"0.01461699455879539"
$ws.Range("F4").Value = [double]"0.00853943451564367"
$ws.Range("G4").Value = [double]"0.01133211216189923"
$ws.Range("H4").Value = [double]"0.09844220277234567"

$ws.Range("B5").Value = [double]"0.05038441735037621"
$ws.Range("C5").Value = [double]"0.001481591126985743"
$ws.Range("D5").Value = [double]"24.0324537832172"
$ws.Range("E5").Value = [double]"0.09869496024145431"
$ws.Range("F5").Value = [double]"0.04748054374149523"
$ws.Range("G5").Value = [double]"0.0532882909592572"
$ws.Range("H5").Value = [double]"0.1388908467839504"

$ws.Range("B6").Value = [double]"0.1435283801441561"
$ws.Range("C6").Value = [double]"0.003027455121449116"
$ws.Range("D6").Value = [double]"38.33826927612975"
$ws.Range("E6").Value = [double]"0.04449515760062214"
$ws.Range("F6").Value = [double]"0.1375946625751373"
$ws.Range("G6").Value = [double]"0.1494620977131749"
$ws.Range("H6").Value = [double]"0.2320348095777303"

$ws.Range("B7").Value = [double]"0.1474985606801619"
$ws.Range("C7").Value = [double]"0.003617413303276194"
$ws.Range("D7").Value = [double]"42.00640300074584"
$ws.Range("E7").Value = [double]"0.01615739860040682"
$ws.Range("F7").Value = [double]"0.1404085446284985"
$ws.Range("G7").Value = [double]"0.1545885767318253"
$ws.Range("H7").Value = [double]"0.2360049901137362"

$ws.Range("B8").Value = [double]"0.1541295751607056"
$ws.Range("C8").Value = [double]"0.006678467174539767"
$ws.Range("D8").Value = [double]"36.87449747958502"
$ws.Range("E8").Value = [double]"0.1611515920560151"
$ws.Range("F8").Value = [double]"0.1410399706319905"
$ws.Range("G8").Value = [double]"0.1672191796894208"
$ws.Range("H8").Value = [double]"0.2426360045942799"

$ws.Range("B9").Value = [double]"0.1497420033581862"
$ws.Range("C9").Value = [double]"-0.0003971853232446947"
$ws.Range("D9").Value = [double]"36.24206586876861"
$ws.Range("E9").Value = [double]"-2.1155163672883"
$ws.Range("F9").Value = [double]"0.1471802994563764"
$ws.Range("G9").Value = [double]"0.1612311264783391"
$ws.Range("H9").Value = [double]"0.2382484327917604"

$ws.Range("B10").Value = [double]"-0.08850642943357422"
$ws.Range("C10").Value = [double]"0.0004388197056097383"
$ws.Range("D10").Value = [double]"-216.0002379275462"
$ws.Range("E10").Value = [double]"0"
$ws.Range("F10").Value = [double]"-0.08936650296156834"
$ws.Range("G10").Value = [double]"-0.08764635590558011"

$ws.Range("B11").Value = [double]"-0.03613024727371009"
$ws.Range("C11").Value = [double]"0.0005018085756394649"
$ws.Range("D11").Value = [double]"-74.94290853472184"
$ws.Range("E11").Value = [double]"0"
$ws.Range("F11").Value = [double]"-0.03711377714342275"
$ws.Range("G11").Value = [double]"-0.03514671740399741"
$ws.Range("H11").Value = [double]"0.05237618215986413"

$ws.Range("B12").Value = [double]"-0.02850736696045795"
$ws.Range("C12").Value = [double]"0.000494573272642717"
$ws.Range("D12").Value = [double]"-61.39002078563549"
$ws.Range("E12").Value = [double]"2.991974925187464e-175"
$ws.Range("F12").Value = [double]"-0.02947671585979815"
$ws.Range("G12").Value = [double]"-0.02753801806111772"
$ws.Range("H12").Value = [double]"0.05999906247311627"

$ws.Range("B13").Value = [double]"-0.02330292712748843"
$ws.Range("C13").Value = [double]"0.0004932577155800657"
$ws.Range("D13").Value = [double]"-50.79828173315393"
$ws.Range("E13").Value = [double]"1.284825300114752e-141"
$ws.Range("F13").Value = [double]"-0.02426969764100124"
$ws.Range("G13").Value = [double]"-0.02233615661397561"
$ws.Range("H13").Value = [double]"0.06520350230608579"

$ws.Range("B14").Value = [double]"-0.01892315855453118"
$ws.Range("C14").Value = [double]"0.0004873558885365723"
$ws.Range("D14").Value = [double]"-41.39414847844651"
$ws.Range("E14").Value = [double]"4.123136802050445e-93"
$ws.Range("F14").Value = [double]"-0.01987836166664428"
$ws.Range("G14").Value = [double]"-0.01796795544241809"
$ws.Range("H14").Value = [double]"0.06958327087904304"

$ws.Range("B15").Value = [double]"-0.01334480337335203"
$ws.Range("C15").Value = [double]"0.0004660309868089372"
$ws.Range("D15").Value = [double]"-31.64008833482075"
$ws.Range("E15").Value = [double]"1.542997232946368e-31"
$ws.Range("F15").Value = [double]"-0.01425821024202413"
$ws.Range("G15").Value = [double]"-0.01243139650467992"
$ws.Range("H15").Value = [double]"0.0751616260602222"

$ws.Range("B16").Value = [double]"-0.01056154666330594"
$ws.Range("C16").Value = [double]"0.0004457147118100469"
$ws.Range("D16").Value = [double]"-25.29313899271546"
$ws.Range("E16").Value = [double]"2.440284701443214e-39"
$ws.Range("F16").Value = [double]"-0.01143513420426613"
$ws.Range("G16").Value = [double]"-0.009687959122345749"
$ws.Range("H16").Value = [double]"0.07794488277026827"

$ws.Range("B17").Value = [double]"-0.008323921296145084"
$ws.Range("C17").Value = [double]"0.0004453057277569389"
$ws.Range("D17").Value = [double]"-20.48058997765746"
$ws.Range("E17").Value = [double]"0.02115549076508775"
$ws.Range("F17").Value = [double]"-0.009196707172414699"
$ws.Range("G17").Value = [double]"-0.007451135419875469"
$ws.Range("H17").Value = [double]"0.08018250813742914"

$ws.Range("B18").Value = [double]"-0.00820476721482506"
$ws.Range("C18").Value = [double]"0.0004772400031403092"
$ws.Range("D18").Value = [double]"-18.01775918915204"
$ws.Range("E18").Value = [double]"9.043266100495062e-16"
$ws.Range("F18").Value = [double]"-0.009140143445066801"
$ws.Range("G18").Value = [double]"-0.00726939098458332"
$ws.Range("H18").Value = [double]"0.08030166221874915"

$ws.Range("B19").Value = [double]"-0.006466934182283023"
$ws.Range("C19").Value = [double]"0.0004753091026563042"
$ws.Range("D19").Value = [double]"-15.60712229935283"
$ws.Range("E19").Value = [double]"0.004876474867258261"
$ws.Range("F19").Value = [double]"-0.007398525975863419"
$ws.Range("G19").Value = [double]"-0.005535342388702629"
$ws.Range("H19").Value = [double]"0.08203949525129119"

$ws.Range("B20").Value = [double]"-0.005543121604593432"
$ws.Range("C20").Value = [double]"0.0004756594907113275"
$ws.Range("D20").Value = [double]"-13.06122375267371"
$ws.Range("E20").Value = [double]"4.171514952288219e-07"
$ws.Range("F20").Value = [double]"-0.006475400112895754"
$ws.Range("G20").Value = [double]"-0.004610843096291113"
$ws.Range("H20").Value = [double]"0.08296330782898079"

$ws.Range("B21").Value = [double]"-0.003263367144238542"
$ws.Range("C21").Value = [double]"0.0004665058733254265"
$ws.Range("D21").Value = [double]"-8.209332715407827"
$ws.Range("E21").Value = [double]"0.0002583707643156257"
$ws.Range("F21").Value = [double]"-0.004177704713698394"
$ws.Range("G21").Value = [double]"-0.002349029574778689"
$ws.Range("H21").Value = [double]"0.08524306228933567"

$ws.Range("B22").Value = [double]"-0.003564553295063925"
$ws.Range("C22").Value = [double]"0.0004792158176284522"
$ws.Range("D22").Value = [double]"-7.675268909836387"
$ws.Range("E22").Value = [double]"0.000714658274040344"
$ws.Range("F22").Value = [double]"-0.004503802103207974"
$ws.Range("G22").Value = [double]"-0.002625304486919877"
$ws.Range("H22").Value = [double]"0.08494187613851029"

$ws.Range("B23").Value = [double]"-0.003557501292084617"
$ws.Range("C23").Value = [double]"0.0004523152716267988"
$ws.Range("D23").Value = [double]"-7.605916162201303"
$ws.Range("E23").Value = [double]"0.02917577412055372"
$ws.Range("F23").Value = [double]"-0.004444025693828793"
$ws.Range("G23").Value = [double]"-0.002670976890340441"
$ws.Range("H23").Value = [double]"0.0849489281414896"

$ws.Range("B24").Value = [double]"-0.002250385406978238"
$ws.Range("C24").Value = [double]"0.0004473134614418279"
$ws.Range("D24").Value = [double]"-5.323828450117258"
$ws.Range("E24").Value = [double]"0.0008650435432193291"
$ws.Range("F24").Value = [double]"-0.003127106371952358"
$ws.Range("G24").Value = [double]"-0.001373664442004118"
$ws.Range("H24").Value = [double]"0.08625604402659598"

$ws.Range("B25").Value = [double]"-0.001708380286252814"
$ws.Range("C25").Value = [double]"0.0004207538742031859"
$ws.Range("D25").Value = [double]"-4.37801352867651"
$ws.Range("E25").Value = [double]"0.1289623274452273"
$ws.Range("F25").Value = [double]"-0.002533045184823293"
$ws.Range("G25").Value = [double]"-0.0008837153876823341"
$ws.Range("H25").Value = [double]"0.08679804914732141"

$ws.Range("B26").Value = [double]"0.2106851903530488"
$ws.Range("C26").Value = [double]"0.003307286072215352"
$ws.Range("D26").Value = [double]"98.88865119218204"
$ws.Range("E26").Value = [double]"0.01755462643176985"
$ws.Range("F26").Value = [double]"0.20420300963343"
$ws.Range("G26").Value = [double]"0.2171673710726675"
$ws.Range("H26").Value = [double]"0.299191619786623"

Write-Host "updated cells"